$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.067.84'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '3.088.91'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.39'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.49'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.086.83'
$ws.Range('E8').Value = '  -1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.44'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.475'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('E13').Value = '  -2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.47'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('E15').Value = '  -1.55%  '
$ws.Range('D16').Value = '3.602.07'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '67.017.63'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.06'
$ws.Range('E18').Value = '  -1.80%  '
$ws.Range('D19').Value = '3.099.05'
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.59'
$ws.Range('E20').Value = '  +2.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '488.77'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.695'
$ws.Range('E22').Value = '  -3.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.76'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.61'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.97'
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.26'
$ws.Range('E26').Value = '  -3.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.41'
$ws.Range('E27').Value = '  +2.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.73'
$ws.Range('E29').Value = '  -3.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.31'
$ws.Range('E30').Value = '  -4.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.65'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.14'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.112'
$ws.Range('E33').Value = '  -2.24%  '
$ws.Range('D34').Value = '0.0₃0942'
$ws.Range('E34').Value = '  -5.86%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.74'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.964'
$ws.Range('E37').Value = '  -2.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '46.50'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.01'
$ws.Range('E39').Value = '  -4.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.123'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.305'
$ws.Range('E41').Value = '  -2.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.42'
$ws.Range('E42').Value = '  -2.96%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.787.68'
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '380.54'
$ws.Range('E44').Value = '  -1.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.54'
$ws.Range('E45').Value = '  -9.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0349'
$ws.Range('E46').Value = '  -2.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '134.83'
$ws.Range('E47').Value = '  -1.28%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.68'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.18'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('E51').Value = '  -2.41%  '
